$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "false start" data rows (old rows 2 and 3), shifting rows 4-5 up
$ws.Rows("2:3").Delete()

# Update selection to reflect the new state (A2:XFD2 and A3:XFD3 multi-selection, active cell A3)
$ws.Range("A3:XFD3,A2:XFD2").Select()
